$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

# Numeric cells
$ws.Cells.Item($row, 1).Value = 131221427          # A16
$ws.Cells.Item($row, 2).Value = 57073               # B16
$ws.Cells.Item($row, 5).Value = 100138              # E16
$ws.Cells.Item($row, 17).Value = 471404             # Q16
$ws.Cells.Item($row, 18).Value = 6543714            # R16
$ws.Cells.Item($row, 19).Value = 20                 # S16

# Plain text cells
$ws.Cells.Item($row, 4).Value = "LC"                # D16
$ws.Cells.Item($row, 6).Value = "Tjäder"            # F16
$ws.Cells.Item($row, 7).Value = "Tetrao urogallus"  # G16
$ws.Cells.Item($row, 8).Value = "Linnaeus, 1758"    # H16
$ws.Cells.Item($row, 13).Value = "färska spår"      # M16
$ws.Cells.Item($row, 16).Value = "SV Rankemossen, Stora Rankemossen, Nrk" # P16
$ws.Cells.Item($row, 20).Value = "Örebro"           # T16
$ws.Cells.Item($row, 21).Value = "Laxå"             # U16
$ws.Cells.Item($row, 22).Value = "Närke"            # V16
$ws.Cells.Item($row, 23).Value = "Skagershult"      # W16
$ws.Cells.Item($row, 26).Value = "17:15"            # Z16
$ws.Cells.Item($row, 28).Value = "17:15"            # AB16
$ws.Cells.Item($row, 29).Value = "Färska spårlöpor i snön."  # AC16
$ws.Cells.Item($row, 49).Value = "Therese Steiner"  # AW16
$ws.Cells.Item($row, 50).Value = "Therese Steiner"  # AX16

# Text cells that would otherwise be misread as numbers/dates -
# force text format first so they stay text like the source file,
# then restore the default (unformatted) style.
$ws.Cells.Item($row, 9).NumberFormat = "@"
$ws.Cells.Item($row, 9).Value = "1"                 # I16
$ws.Cells.Item($row, 9).Style = "Normal"

$ws.Cells.Item($row, 25).NumberFormat = "@"
$ws.Cells.Item($row, 25).Value = "2026-02-18"       # Y16
$ws.Cells.Item($row, 25).Style = "Normal"

$ws.Cells.Item($row, 27).NumberFormat = "@"
$ws.Cells.Item($row, 27).Value = "2026-02-18"       # AA16
$ws.Cells.Item($row, 27).Style = "Normal"

# Boolean cells
$ws.Cells.Item($row, 30).Value = $false             # AD16
$ws.Cells.Item($row, 31).Value = $false             # AE16
$ws.Cells.Item($row, 33).Value = $false             # AG16

# Empty-string cells (AT16/AY16 are blank in the source row; re-apply
# the default style so the cell is still materialised in the sheet,
# matching the self-closed empty cell the source file has for this
# column on every other row).
$ws.Cells.Item($row, 46).Value = ""                 # AT16
$ws.Cells.Item($row, 46).Style = "Normal"
$ws.Cells.Item($row, 51).Value = ""                 # AY16
$ws.Cells.Item($row, 51).Style = "Normal"
